$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Set header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Set data row
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Set column B width to fit content
$ws.Columns.Item(2).ColumnWidth = 11.140625

# Zoom and selection
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 250
